{"js": "// Update the date label and each multiplication problem's operands.\nconst replacements = [\n  [\"2025-10-19 Sunday\", \"2025-10-20 Monday\"],\n  [\"445\u00d78=\", \"637\u00d72=\"],\n  [\"238\u00d78=\", \"947\u00d73=\"],\n  [\"880\u00d73=\", \"266\u00d73=\"],\n  [\"901\u00d75=\", \"203\u00d73=\"],\n  [\"297\u00d73=\", \"710\u00d76=\"],\n  [\"507\u00d77=\", \"375\u00d77=\"],\n  [\"842\u00d76=\", \"159\u00d77=\"],\n  [\"694\u00d77=\", \"145\u00d75=\"],\n  [\"289\u00d77=\", \"928\u00d79=\"],\n  [\"331\u00d79=\", \"891\u00d78=\"],\n  [\"804\u00d73=\", \"141\u00d75=\"],\n  [\"287\u00d73=\", \"510\u00d74=\"],\n  [\"569\u00d79=\", \"546\u00d75=\"],\n  [\"955\u00d75=\", \"683\u00d79=\"],\n  [\"558\u00d75=\", \"521\u00d72=\"],\n  [\"324\u00d74=\", \"566\u00d74=\"],\n  [\"276\u00d78=\", \"186\u00d74=\"],\n  [\"448\u00d75=\", \"439\u00d77=\"],\n  [\"875\u00d77=\", \"254\u00d75=\"],\n  [\"395\u00d79=\", \"989\u00d73=\"],\n  [\"411\u00d76=\", \"848\u00d79=\"],\n  [\"847\u00d76=\", \"649\u00d79=\"],\n  [\"797\u00d74=\", \"705\u00d79=\"],\n  [\"148\u00d72=\", \"144\u00d79=\"],\n  [\"962\u00d72=\", \"918\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and each multiplication problem's operands.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-19 Sunday\", \"2025-10-20 Monday\"),\n    @(\"445\u00d78=\", \"637\u00d72=\"),\n    @(\"238\u00d78=\", \"947\u00d73=\"),\n    @(\"880\u00d73=\", \"266\u00d73=\"),\n    @(\"901\u00d75=\", \"203\u00d73=\"),\n    @(\"297\u00d73=\", \"710\u00d76=\"),\n    @(\"507\u00d77=\", \"375\u00d77=\"),\n    @(\"842\u00d76=\", \"159\u00d77=\"),\n    @(\"694\u00d77=\", \"145\u00d75=\"),\n    @(\"289\u00d77=\", \"928\u00d79=\"),\n    @(\"331\u00d79=\", \"891\u00d78=\"),\n    @(\"804\u00d73=\", \"141\u00d75=\"),\n    @(\"287\u00d73=\", \"510\u00d74=\"),\n    @(\"569\u00d79=\", \"546\u00d75=\"),\n    @(\"955\u00d75=\", \"683\u00d79=\"),\n    @(\"558\u00d75=\", \"521\u00d72=\"),\n    @(\"324\u00d74=\", \"566\u00d74=\"),\n    @(\"276\u00d78=\", \"186\u00d74=\"),\n    @(\"448\u00d75=\", \"439\u00d77=\"),\n    @(\"875\u00d77=\", \"254\u00d75=\"),\n    @(\"395\u00d79=\", \"989\u00d73=\"),\n    @(\"411\u00d76=\", \"848\u00d79=\"),\n    @(\"847\u00d76=\", \"649\u00d79=\"),\n    @(\"797\u00d74=\", \"705\u00d79=\"),\n    @(\"148\u00d72=\", \"144\u00d79=\"),\n    @(\"962\u00d72=\", \"918\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
